$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing account 001922009 / SOFIA / 10514.44 (row 4: header is row 1)
$ws.Rows.Item(4).Delete()
